$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition the workbook window (matches the saved bookViews entry) ---
$win = $wb.Windows.Item(1)
$win.Left = 11715
$win.Top = 3420

# --- Update the JSSP problem data (processing order / times matrix) ---
$ws.Cells.Item(1,1).Value = 1
$ws.Cells.Item(1,2).Value = 2
$ws.Cells.Item(1,3).Value = 0
$ws.Cells.Item(1,4).Value = 3
$ws.Cells.Item(1,5).Value = 2
$ws.Cells.Item(1,6).Value = 2

$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = 1
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 2
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 4

$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 3
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = 4
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0

# --- Extend the highlighted/formatted block that used to be just H6
#     down into a big A4:L20 grid, re-using the same cell style (format painter) ---
$ws.Range("H6").Copy()
$ws.Range("A4:L20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the saved selection/active cell ---
$ws.Range("K26").Select() | Out-Null
